$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.362.85'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').Value = '3.381.72'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.51'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.16%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.381.81'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.47'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.123'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.387'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').Value = '3.958.69'
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.98'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '3.382.00'
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('E17').Value = '  -3.28%  '
$ws.Range('D18').Value = '60.460.62'
$ws.Range('E18').Value = '  -1.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.85'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '386.66'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.555'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.22'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -6.82%  '
$ws.Range('D27').Value = '3.529.22'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.178'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  -5.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.91'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('E32').Value = '  -2.16%  '
$ws.Range('E33').Value = '  -8.55%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.56'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('D36').Value = '3.412.11'
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '167.53'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.41%  '
$ws.Range('E40').Value = '  -6.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0769'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.82'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.780'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.42'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.95%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.23'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.53%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.68'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.84%  '
$ws.Range('D48').Value = '2.524.10'
$ws.Range('E48').Value = '  -3.10%  '
$ws.Range('E49').Value = '  -4.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.17'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  -3.92%  '
